$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 301.625
$ws.Range("J9").Value = 361.5
$ws.Range("L9").Value = 361.5
$ws.Range("N9").Value = -699.5
$ws.Range("H62").Value = 3484.5
$ws.Range("I62").Value = 3484.5
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3484.5
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -2860.5
$ws.Range("H65").Value = 3484.5
$ws.Range("I65").Value = 3484.5
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 17422.5
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -14302.5
$ws.Range("H118").Value = 299.25
$ws.Range("I118").Value = 299.66666
$ws.Range("J118").Value = 298
$ws.Range("K118").Value = 898.9999799999999
$ws.Range("L118").Value = 894
$ws.Range("M118").Value = 758.0000200000001
$ws.Range("N118").Value = -4208
$ws.Range("H132").Value = 1328.4584
$ws.Range("I132").Value = 1294.4762
$ws.Range("K132").Value = 3883.4286
$ws.Range("M132").Value = -1353.4286
$ws.Range("H138").Value = 4745.9536
$ws.Range("I138").Value = 3175.3235
$ws.Range("J138").Value = 6468.5806
$ws.Range("K138").Value = 9525.970499999999
$ws.Range("L138").Value = 19405.7418
$ws.Range("M138").Value = -4385.970499999999
$ws.Range("N138").Value = -29685.7418

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3874.1428
$ws.Range("I2").Value = 2511
$ws.Range("K2").Value = 2511
$ws.Range("M2").Value = -2398
$ws.Range("H26").Value = 2756.3076
$ws.Range("I26").Value = 1313.6666
$ws.Range("J26").Value = 3992.8572
$ws.Range("K26").Value = 1313.6666
$ws.Range("L26").Value = 3992.8572
$ws.Range("M26").Value = -983.6666
$ws.Range("N26").Value = -4652.8572
$ws.Range("H32").Value = 8779.591
$ws.Range("I32").Value = 6008.6113
$ws.Range("K32").Value = 6008.6113
$ws.Range("M32").Value = -5721.6113
$ws.Range("H35").Value = 1478.1428
$ws.Range("I35").Value = 1532
$ws.Range("J35").Value = 1155
$ws.Range("K35").Value = 1532
$ws.Range("L35").Value = 1155
$ws.Range("M35").Value = -1126
$ws.Range("N35").Value = -1967
$ws.Range("H74").Value = 1200
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 1200
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H102").Value = 822.2
$ws.Range("I102").Value = 822.2
$ws.Range("K102").Value = 822.2
$ws.Range("M102").Value = 799.8
$ws.Range("H116").Value = 3874.1428
$ws.Range("I116").Value = 2511
$ws.Range("K116").Value = 2511
$ws.Range("M116").Value = -217
$ws.Range("H132").Value = 1327.2059
$ws.Range("I132").Value = 1161
$ws.Range("J132").Value = 1788.8889
$ws.Range("K132").Value = 3483
$ws.Range("L132").Value = 5366.6667
$ws.Range("M132").Value = -953
$ws.Range("N132").Value = -10426.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3874.1428
$ws.Range("I3").Value = 2511
$ws.Range("K3").Value = 2511
$ws.Range("M3").Value = -2397
$ws.Range("H86").Value = 3362.7334
$ws.Range("I86").Value = 2828.4167
$ws.Range("K86").Value = 2828.4167
$ws.Range("M86").Value = -1705.4167
$ws.Range("H89").Value = 3362.7334
$ws.Range("I89").Value = 2828.4167
$ws.Range("K89").Value = 14142.0835
$ws.Range("M89").Value = -8526.083500000001
$ws.Range("H99").Value = 54157.74
$ws.Range("I99").Value = 64032.375
$ws.Range("K99").Value = 64032.375
$ws.Range("M99").Value = -62534.375
$ws.Range("H105").Value = 4116
$ws.Range("I105").Value = 4116
$ws.Range("K105").Value = 4116
$ws.Range("M105").Value = -2369
$ws.Range("H107").Value = 1436.6666
$ws.Range("I107").Value = 1436.6666
$ws.Range("K107").Value = 1436.6666
$ws.Range("M107").Value = 483.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 208.23529
$ws.Range("I7").Value = 155.46153
$ws.Range("J7").Value = 379.75
$ws.Range("K7").Value = 155.46153
$ws.Range("L7").Value = 379.75
$ws.Range("M7").Value = -42.46153000000001
$ws.Range("N7").Value = -605.75
$ws.Range("H22").Value = 200
$ws.Range("I22").Value = 200
$ws.Range("K22").Value = 200
$ws.Range("M22").Value = 150
$ws.Range("H58").Value = 2257.0625
$ws.Range("I58").Value = 1300
$ws.Range("K58").Value = 1300
$ws.Range("M58").Value = -1097
$ws.Range("H132").Value = 1488.2
$ws.Range("I132").Value = 1480.3334
$ws.Range("K132").Value = 4441.0002
$ws.Range("M132").Value = -1911.0002
$ws.Range("H134").Value = 1846.8182
$ws.Range("I134").Value = 1445.875
$ws.Range("J134").Value = 2916
$ws.Range("K134").Value = 4337.625
$ws.Range("L134").Value = 8748
$ws.Range("M134").Value = -1802.625
$ws.Range("N134").Value = -13818
$ws.Range("H136").Value = 2257.0625
$ws.Range("I136").Value = 1300
$ws.Range("K136").Value = 3900
$ws.Range("M136").Value = -1350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 81.666664
$ws.Range("I17").Value = 105
$ws.Range("K17").Value = 315
$ws.Range("M17").Value = -146
$ws.Range("H56").Value = 12287.182
$ws.Range("I56").Value = 12287.182
$ws.Range("K56").Value = 12287.182
$ws.Range("M56").Value = -11757.182
$ws.Range("H75").Value = 425
$ws.Range("I75").Value = 425
$ws.Range("K75").Value = 1275
$ws.Range("M75").Value = -277
$ws.Range("H78").Value = 425
$ws.Range("I78").Value = 425
$ws.Range("K78").Value = 3825
$ws.Range("M78").Value = 1167
$ws.Range("H116").Value = 1200
$ws.Range("I116").Value = 1200
$ws.Range("K116").Value = 3600
$ws.Range("M116").Value = -158
$ws.Range("H132").Value = 5156.1763
$ws.Range("I132").Value = 5212.4287
$ws.Range("J132").Value = 4893.6665
$ws.Range("K132").Value = 46911.85830000001
$ws.Range("L132").Value = 44042.9985
$ws.Range("M132").Value = -44381.85830000001
$ws.Range("N132").Value = -49102.9985
$ws.Range("H137").Value = 3521.1667
$ws.Range("I137").Value = 3114
$ws.Range("J137").Value = 8000
$ws.Range("K137").Value = 9342
$ws.Range("L137").Value = 24000
$ws.Range("M137").Value = -4242
$ws.Range("N137").Value = -34200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 656.1667
$ws.Range("I2").Value = 46.4
$ws.Range("K2").Value = 46.4
$ws.Range("M2").Value = 66.59999999999999
$ws.Range("H70").Value = 7621
$ws.Range("I70").Value = 7621
$ws.Range("K70").Value = 7621
$ws.Range("M70").Value = -7351
$ws.Range("H73").Value = 7621
$ws.Range("I73").Value = 7621
$ws.Range("K73").Value = 7621
$ws.Range("M73").Value = -6685
$ws.Range("H132").Value = 2699.4348
$ws.Range("I132").Value = 2232.9333
$ws.Range("J132").Value = 3574.125
$ws.Range("K132").Value = 6698.7999
$ws.Range("L132").Value = 10722.375
$ws.Range("M132").Value = -4168.7999
$ws.Range("N132").Value = -15782.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 12502.5
$ws.Range("I18").Value = 12502.5
$ws.Range("K18").Value = 12502.5
$ws.Range("M18").Value = -12330.5
$ws.Range("H30").Value = 3501.25
$ws.Range("I30").Value = 3501.25
$ws.Range("K30").Value = 3501.25
$ws.Range("M30").Value = -3393.25
$ws.Range("H64").Value = 23750
$ws.Range("J64").Value = 23750
$ws.Range("L64").Value = 23750
$ws.Range("N64").Value = -24200
$ws.Range("H67").Value = 23750
$ws.Range("J67").Value = 23750
$ws.Range("L67").Value = 23750
$ws.Range("N67").Value = -25310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32023.334
$ws.Range("J54").Value = 23000
$ws.Range("L54").Value = 23000
$ws.Range("N54").Value = -24040
$ws.Range("H68").Value = 31249.5
$ws.Range("J68").Value = 31249.5
$ws.Range("L68").Value = 31249.5
$ws.Range("N68").Value = -32871.5
$ws.Range("H71").Value = 31249.5
$ws.Range("J71").Value = 31249.5
$ws.Range("L71").Value = 93748.5
$ws.Range("N71").Value = -101860.5
$ws.Range("H81").Value = 7262.304
$ws.Range("I81").Value = 1536
$ws.Range("K81").Value = 3072
$ws.Range("M81").Value = -2011
$ws.Range("H84").Value = 7262.304
$ws.Range("I84").Value = 1536
$ws.Range("K84").Value = 15360
$ws.Range("M84").Value = -10056
$ws.Range("H126").Value = 2158.2273
$ws.Range("I126").Value = 1478.4286
$ws.Range("K126").Value = 4435.2858
$ws.Range("M126").Value = -1965.2858
